$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) columns to snake_case names
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

$textInfo = (Get-Culture).TextInfo

$oldPct = 0.0009708737864077669
$newPct = 0.0009708737864077668

# 2. Walk the data rows (2 through 573) and:
#    - Title-case the place names in columns A and B
#    - Fix the 3/3090 percentage rounding value in column D
for ($r = 2; $r -le 573; $r++) {
    $colA = $ws.Cells.Item($r, 1)
    $valA = $colA.Value()
    if ($valA -ne $null -and $valA -is [string]) {
        $colA.Value = $textInfo.ToTitleCase($valA)
    }

    $colB = $ws.Cells.Item($r, 2)
    $valB = $colB.Value()
    if ($valB -ne $null -and $valB -is [string]) {
        $colB.Value = $textInfo.ToTitleCase($valB)
    }

    $colD = $ws.Cells.Item($r, 4)
    $valD = $colD.Value()
    if ($valD -ne $null -and $valD -is [double]) {
        if ($valD -eq $oldPct) {
            $colD.Value = $newPct
        }
    }
}

# 3. Remove the trailing footnote rows (575-579), which are no longer part
#    of the cleaned table.
$ws.Range("A575:A579").EntireRow.Delete()
